$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the bracket "winner" cells (round 1 -> round 2 picks) that were
# left as placeholder "p" text - pick the winning team name for each pair.
$ws.Range("E2").Value  = "Pride of Lionesses"
$ws.Range("E4").Value  = "Sneak of Weasels"
$ws.Range("E6").Value  = "Conspiracy of Lemurs"
$ws.Range("E8").Value  = "Lodge of Beavers"
$ws.Range("E10").Value = "Wisdom of Wombats"
$ws.Range("E12").Value = "Cauldron of Bats"
$ws.Range("E14").Value = "Stench of Skunks"
$ws.Range("E16").Value = "Embarrassment of Pandas"

# Restore the view: scroll so row 4 is the top visible row, and leave the
# selection on H15 (matches the saved window state in the workbook).
$ws.Activate()
$ws.Range("H15").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
